$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Replace-InParaLast($paraIndex, $old, $new) {
    # Replace the LAST occurrence of $old within the given paragraph's text with $new.
    $p = $d.Paragraphs.Item($paraIndex)
    $full = $p.Range.Text
    $idx = $full.LastIndexOf($old)
    if ($idx -lt 0) {
        return
    }
    $start = $p.Range.Start + $idx
    $end = $start + $old.Length
    $r = $d.Range($start, $end)
    $r.Text = $new
}

# ---------------------------------------------------------------------------
# Title
Replace-Text "Quantum Revelations: Unraveling Nature's Enigmatic Secrets" "Delving into the World of Politics: A Comprehensive Guide"

# Author
Replace-Text "Oliver Saunders" "Ethan Bennett"

# Email paragraph (paragraph 3): oliver | . | saunders96@emailcentral | . | net
# Target: bennetthan@gnail | . | cam   (last "." run and "net" run are deleted)
Replace-InParaLast 3 "net" ""
Replace-InParaLast 3 "." ""
Replace-InParaLast 3 "saunders96@emailcentral" "cam"
Replace-InParaLast 3 "oliver" "bennetthan@gnail"

# ---------------------------------------------------------------------------
# Body paragraph (paragraph 5) - sub-zone 1 (before first <br/>)
Replace-Text "As we embark on a captivating odyssey into the realm of quantum mechanics, a thought-provoking dance of particles and waves, the very foundations of our universe are unveiled" "Politics, a multifaceted and ever-evolving sphere of human interaction, holds immense significance in shaping our world"
Replace-Text " The subatomic realm, an enigmatic tapestry woven with probability and uncertainty, beckons us to explore phenomena that defy classical intuition" " It encompasses the systems, institutions, and processes through which societies are governed, decisions are made, and resources are allocated"
Replace-Text " From the intrinsic interconnectedness of entangled particles to the tunnel-like passage of particles through impassable barriers, quantum mechanics has revolutionized our understanding of the universe and profoundly influenced diverse fields, ranging from computing to cryptography" " Understanding politics is crucial for young minds as it equips them with the knowledge and skills necessary to navigate the intricacies of a complex and interconnected world"

# sub-zone 2 (between the two <br/><br/> pairs)
Replace-Text "In the quantum domain, the classical laws governing the macroscopic world unravel, leaving us with a fascinating enigma" "In this essay, we will embark on a journey to explore the multifaceted nature of politics"
Replace-Text " Particles exhibit paradoxical behaviors, simultaneously existing and occupying multiple states, an ethereal symphony of possibilities" " We will delve into its fundamental principles, analyze different political ideologies, and examine the diverse structures and institutions that govern nations"
# Remove the trailing two sentences (". The tension...determinism. The quantum realm...understanding") but keep the final "."
Replace-Text ". The tension between competing interpretations gives rise to lively debates and captivating paradoxes, challenging our conventional notions of time, locality, and determinism. The quantum realm beckons us to confront these paradoxical intricacies, pushing the boundaries of human knowledge and understanding" ""
Replace-Text " These mind-bending phenomena have led to the formulation of profound interpretations, questioning the fabric of reality itself" " By gaining a deeper understanding of politics, we can foster active and informed citizenship, enabling individuals to participate effectively in shaping their communities and societies"

# sub-zone 3 (after the second <br/><br/>)
Replace-Text "Quantum mechanics has propelled advancements across disciplines" "The study of politics begins with an exploration of its core concepts"
Replace-Text " In the realm of cryptography, it holds the promise of unbreakable encryption, as the inherent uncertainty of quantum systems renders eavesdropping impossible" " We will investigate the nature of power, authority, and legitimacy, delving into the various ways in which individuals and institutions exercise influence over others"
Replace-Text " This transformative potential has inspired the creation of quantum communication networks, poised to revolutionize secure information transfer" " We will examine the relationship between politics and economics, considering how economic structures and policies impact political decision-making"
# Remove the trailing sentence (". These transformative applications...achievement") but keep the final "."
Replace-Text ". These transformative applications are poised to usher in an era of unprecedented technological progress and redefine the boundaries of human achievement" ""
Replace-Text " Furthermore, quantum computing, harnessing the power of quantum entanglement and superposition, offers the prospect of exponential gains in computational capacity, opening uncharted frontiers for simulations and optimizations" " Additionally, we will analyze the role of political culture and ideology in shaping the political landscape of societies"

# ---------------------------------------------------------------------------
# "Summary" heading - force a real replace to drop <w:lastRenderedPageBreak/>
Replace-Text "Summary" "Summary"

# ---------------------------------------------------------------------------
# Summary body paragraph (paragraph 7)
Replace-Text "In this exploration of quantum mechanics, we have delved into the perplexing realm of subatomic particles, uncovering phenomena that challenge conventional wisdom" "This essay provides a comprehensive overview of politics, exploring its fundamental principles, ideologies, and structures"
Replace-Text " Quantum mechanics has redefined our understanding of nature, fostering lively debates about reality's fundamental structure" " It highlights the importance of political engagement and emphasizes the role of individuals in shaping political outcomes"
# Remove the trailing sentence (". As we continue to unravel...tapestry") but keep the final "."
Replace-Text ". As we continue to unravel the enigmatic secrets of the quantum realm, we may one day resolve its paradoxes and forge a deeper connection with the universe's enigmatic tapestry" ""
Replace-Text " It has fueled remarkable innovations in computing, cryptography, and other fields, hinting at a future brimming with possibilities" " By understanding the complexities of politics, young minds can become active and informed citizens, capable of making informed decisions and contributing positively to their communities and societies"

# ---------------------------------------------------------------------------
# Add an empty paragraph after the Summary body paragraph (before the final section break)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
